$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names / links) that don't look numeric: safe to assign directly ---
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

# --- Numeric-looking text cells (price / volume% / hour) ---
# These must stay TEXT (not be reinterpreted as numbers) to preserve exact
# formatting (trailing zeros, percent signs, etc.), matching the source data's
# inline-string cell type. Pre-formatting the cell as Text ('@') before the
# assignment keeps Excel from silently converting the literal to a Double.
$textCells = @(
    "D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5", "E5", "G5", "D6", "E6", "G6",
    "D7", "E7", "G7", "D8", "E8", "G8", "D9", "E9", "G9", "G10", "D11", "E11", "G11", "D12",
    "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14", "D15", "E15", "G15", "D16", "E16",
    "G16", "D17", "E17", "G17", "D18", "E18", "G18", "D19", "E19", "G19", "D20", "E20", "G20",
    "D21", "E21", "G21", "D22", "E22", "G22", "D23", "E23", "G23", "D24", "E24", "G24", "D25",
    "E25", "G25", "D26", "E26", "G26", "E27", "G27", "D28", "E28", "G28", "G29", "G30", "G31",
    "G32", "G33", "G34", "G35", "G36", "G37", "G38", "G39", "D40", "E40", "G40", "D41", "E41",
    "G41", "D42", "E42", "G42", "D43", "E43", "G43", "D44", "E44", "G44", "D45", "E45", "G45",
    "E46", "G46", "D47", "E47", "G47", "D48", "E48", "G48", "D49", "E49", "G49", "D50", "E50",
    "G50", "G51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '245.10'
$ws.Range("E2").Value = '-0.57%'
$ws.Range("G2").Value = '18'
$ws.Range("D3").Value = '29.04'
$ws.Range("E3").Value = '-2.68%'
$ws.Range("G3").Value = '18'
$ws.Range("D4").Value = '5.241'
$ws.Range("E4").Value = '1.45%'
$ws.Range("G4").Value = '18'
$ws.Range("D5").Value = '0.05699'
$ws.Range("E5").Value = '-0.30%'
$ws.Range("G5").Value = '18'
$ws.Range("D6").Value = '6.617'
$ws.Range("E6").Value = '0.25%'
$ws.Range("G6").Value = '18'
$ws.Range("D7").Value = '3.193'
$ws.Range("E7").Value = '3.60%'
$ws.Range("G7").Value = '18'
$ws.Range("D8").Value = '0.8501'
$ws.Range("E8").Value = '-0.83%'
$ws.Range("G8").Value = '18'
$ws.Range("D9").Value = '0.8536'
$ws.Range("E9").Value = '-1.76%'
$ws.Range("G9").Value = '18'
$ws.Range("G10").Value = '18'
$ws.Range("D11").Value = '0.07072'
$ws.Range("E11").Value = '0.08%'
$ws.Range("G11").Value = '18'
$ws.Range("D12").Value = '0.03156'
$ws.Range("E12").Value = '8.03%'
$ws.Range("G12").Value = '18'
$ws.Range("D13").Value = '0.09216'
$ws.Range("E13").Value = '-1.86%'
$ws.Range("G13").Value = '18'
$ws.Range("D14").Value = '0.001541'
$ws.Range("E14").Value = '1.26%'
$ws.Range("G14").Value = '18'
$ws.Range("D15").Value = '0.0005939'
$ws.Range("E15").Value = '-94.24%'
$ws.Range("G15").Value = '18'
$ws.Range("D16").Value = '0.005891'
$ws.Range("E16").Value = '-2.28%'
$ws.Range("G16").Value = '18'
$ws.Range("D17").Value = '3.493'
$ws.Range("E17").Value = '0.15%'
$ws.Range("G17").Value = '18'
$ws.Range("D18").Value = '2.175'
$ws.Range("E18").Value = '-4.41%'
$ws.Range("G18").Value = '18'
$ws.Range("D19").Value = '0.3164'
$ws.Range("E19").Value = '-0.31%'
$ws.Range("G19").Value = '18'
$ws.Range("D20").Value = '0.03214'
$ws.Range("E20").Value = '-3.91%'
$ws.Range("G20").Value = '18'
$ws.Range("D21").Value = '0.1276'
$ws.Range("E21").Value = '-1.84%'
$ws.Range("G21").Value = '18'
$ws.Range("D22").Value = '3.524'
$ws.Range("E22").Value = '1.71%'
$ws.Range("G22").Value = '18'
$ws.Range("D23").Value = '0.04083'
$ws.Range("E23").Value = '-2.21%'
$ws.Range("G23").Value = '18'
$ws.Range("D24").Value = '0.1379'
$ws.Range("E24").Value = '-0.10%'
$ws.Range("G24").Value = '18'
$ws.Range("D25").Value = '0.001220'
$ws.Range("E25").Value = '-0.01%'
$ws.Range("G25").Value = '18'
$ws.Range("D26").Value = '0.004140'
$ws.Range("E26").Value = '-17.60%'
$ws.Range("G26").Value = '18'
$ws.Range("E27").Value = '-0.83%'
$ws.Range("G27").Value = '18'
$ws.Range("D28").Value = '0.0001449'
$ws.Range("E28").Value = '-98.07%'
$ws.Range("G28").Value = '18'
$ws.Range("G29").Value = '18'
$ws.Range("G30").Value = '18'
$ws.Range("G31").Value = '18'
$ws.Range("G32").Value = '18'
$ws.Range("G33").Value = '18'
$ws.Range("G34").Value = '18'
$ws.Range("G35").Value = '18'
$ws.Range("G36").Value = '18'
$ws.Range("G37").Value = '18'
$ws.Range("G38").Value = '18'
$ws.Range("G39").Value = '18'
$ws.Range("D40").Value = '0.03756'
$ws.Range("E40").Value = '0.27%'
$ws.Range("G40").Value = '18'
$ws.Range("D41").Value = '0.1064'
$ws.Range("E41").Value = '-0.58%'
$ws.Range("G41").Value = '18'
$ws.Range("D42").Value = '0.003708'
$ws.Range("E42").Value = '-35.96%'
$ws.Range("G42").Value = '18'
$ws.Range("D43").Value = '0.002490'
$ws.Range("E43").Value = '24.51%'
$ws.Range("G43").Value = '18'
$ws.Range("D44").Value = '0.009356'
$ws.Range("E44").Value = '12.68%'
$ws.Range("G44").Value = '18'
$ws.Range("D45").Value = '0.00005271'
$ws.Range("E45").Value = '1.16%'
$ws.Range("G45").Value = '18'
$ws.Range("E46").Value = '-0.03%'
$ws.Range("G46").Value = '18'
$ws.Range("D47").Value = '0.07498'
$ws.Range("E47").Value = '29.27%'
$ws.Range("G47").Value = '18'
$ws.Range("D48").Value = '0.002439'
$ws.Range("E48").Value = '-5.20%'
$ws.Range("G48").Value = '18'
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").Value = '-0.03%'
$ws.Range("G49").Value = '18'
$ws.Range("D50").Value = '0.0002000'
$ws.Range("E50").Value = '-0.03%'
$ws.Range("G50").Value = '18'
$ws.Range("G51").Value = '18'
